$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "CCACTCAACCATCCACTCCC"
$ws.Range("B3").Value = "AAGGTGAGTGTGGCAAGTGG"
$ws.Range("B4").Value = "ACCCACACACACACAACACT"
$ws.Range("B5").Value = "CAGGGTAAGTGGCAGTGGAG"
$ws.Range("B6").Value = "TCACTCTCCAACTTCTCTGCT"
$ws.Range("B7").Value = "AGCAGATTTCGGAGGTGTGG"
